# Generate Report for Handoff
# A new handoff/xliff-generation pass completed for the
# "892fead8-f8cd-4564-979c-3a1e38d4a1fd.md" file (row 6 in each sheet),
# so its recorded timestamps move forward:
#   - Overview!G6  (Latest HO Xliff Generate Date)  -> 2016-09-05 12:52:22
#   - zh-cn!H6     (Latest Handoff Datetime)         -> 2016-09-05 12:52:17
#   - de-de!H6     (Latest Handoff Datetime)         -> 2016-09-05 12:52:22

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G6").Value = "2016-09-05 12:52:22"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H6").Value = "2016-09-05 12:52:17"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H6").Value = "2016-09-05 12:52:22"
